$d = $word.ActiveDocument
$RQ = [char]0x2019   # right single quotation mark (')

# ===========================================================================
# STEP 1 -- text-only edits (do not change paragraph count / ordering)
# ===========================================================================

# --- PL-4 a) CivicActions --------------------------------------------------
$d.Content.Find.Execute(
    "an Acceptable Use Policy that describes their responsibilities and expected behavior with regard to information and information system usage. This information is available in the CivicActions Handbook under Security Policy here:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "the rules that describes their responsibilities and expected behavior with regard to information and information system usage. These rules, defined as the Acceptable Use Policy, are included in the CivicActions Security Policy accessible here :",
    2) | Out-Null

$d.Content.Find.Execute(
    "and has also been uploaded to CSAM as",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "which has also been uploaded to CSAM as",
    2) | Out-Null

$d.Content.Find.Execute(
    "Appendix J1 - CivicActions Security Policy 20190226.docx",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Appendix J1 - System Rules of Behavior - Privileged User",
    2) | Out-Null

$d.Content.Find.Execute(
    "Privileged User" + $RQ + ".",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Privileged User" + $RQ + " (CivicActions Security Policy 20190226.docx).",
    2) | Out-Null

# --- PL-4 b) CivicActions ---------------------------------------------------
$pb = $d.Paragraphs(40).Range.Duplicate
$pb.Find.Execute(
    "All CivicActions employees are required to read and sign the Security Policy*as artifact:",
    $true, $false, $true, $false, $false, $true, 1, $false,
    "CivicActions HR receives a signed acknowledgment from all employees, indicating that they have read, understand, and agree to abide by the rules of behavior, before authorizing access to information and the information system. The text of the electronically signed (via DocuSign) acknowledgement document has been uploaded to CSAM as artifact:",
    2) | Out-Null

# --- PL-4 c) CivicActions ---------------------------------------------------
$d.Content.Find.Execute(
    "The CivicActions Acceptable Use Policy/Rules of Behavior are reviewed by CivicActions Security and Operations at least annually and is updated at least every three years.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "CivicActions reviews the CivicActions Security Policy at least annually and updates is as required.",
    2) | Out-Null

# --- PL-4 d) CivicActions ---------------------------------------------------
$d.Content.Find.Execute(
    "CivicActions employees re-sign acknowledgement of the Acceptable Use/Rules of Behavior policy document whenever significant changes are made. The Director of Human Resources retains the signed acknowledgements.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "CivicActions requires individuals who have signed a previous version of the CivicActions Security Policy to read and re-sign when any part of it, including the Acceptable Use Policy/Rules of Behavior, are revised/updated.",
    2) | Out-Null

Write-Host "step 1 complete"

# ===========================================================================
# STEP 2 -- insert the new "LINCS" Heading5 + FirstParagraph blocks.
# Processed bottom-to-top so paragraph indices used below stay valid.
# ===========================================================================

function Add-LincsBlock($afterParaIndex, $bookmarkName, $bodyText) {

    $anchor = $d.Paragraphs($afterParaIndex).Range
    $anchor.InsertParagraphAfter() | Out-Null

    $headingIndex = $afterParaIndex + 1
    $bodyIndex = $afterParaIndex + 2

    $d.Paragraphs($headingIndex).Range.InsertParagraphAfter() | Out-Null

    $h = $d.Paragraphs($headingIndex)
    $h.Style = "Heading 5"
    $h.Range.Text = "LINCS"

    $b = $d.Paragraphs($bodyIndex)
    $b.Style = "First Paragraph"
    $b.Range.Text = $bodyText

    $bmRange = $d.Paragraphs($headingIndex).Range
    $bmRange.MoveEnd(1, -1) | Out-Null
    $d.Bookmarks.Add($bookmarkName, $bmRange) | Out-Null
}

# after d) CivicActions (para 46) -- end of PL-4 section
Add-LincsBlock 46 "lincs-5" "LINCS requires individuals who have signed a previous version of the rules of behavior to read and re-sign when the Rules of Behavior are revised/updated."

# after c) CivicActions (para 43)
Add-LincsBlock 43 "lincs-4" "LINCS reviews the Rules of Behavior at least annually and updates it as required."

# after b) CivicActions (para 40)
Add-LincsBlock 40 "lincs-3" "The LINCS System Owner receives a signed acknowledgment from such individuals that are not CivicActions employees, indicating that they have read, understand, and agree to abide by the rules of behavior, before authorizing access to information and the information system."

# after a) CivicActions (para 37)
$lincs2Body = "LINCS has created and made readily available to individuals requiring access to the information system the rules that describes their responsibilities and expected behavior with regard to information and information system usage. These rules are captured in " + [char]0x2018 + "Appendix J2 - System Rules of Behavior - General User" + $RQ + " (LINCSSystemRoB2019-template.docx). LINCS has reviewed and accepted as a superset alternative the CivicActions Acceptable Use Policy."
Add-LincsBlock 37 "lincs-2" $lincs2Body

Write-Host "step 2 (four LINCS blocks) complete"
